$wb = $excel.ActiveWorkbook

# --- Sheet "config": update file paths ---
$cfg = $wb.Worksheets.Item("config")

$cfg.Range("B2").Value = "/home/ohel/Schreibtisch/FAST_sim/sim/5MW_Land_DLL_WTurb.fst"
$cfg.Range("B3").Value = "/home/ohel/Schreibtisch/FAST_sim/sim/NRELOffshrBsline5MW_Onshore_ElastoDyn.dat"
$cfg.Range("B4").Value = "/home/ohel/Schreibtisch/FAST_sim/sim/NRELOffshrBsline5MW_Onshore_ServoDyn.dat"
$cfg.Range("B5").Value = "/home/ohel/Schreibtisch/FAST_sim/sim/NRELOffshrBsline5MW_Onshore_AeroDyn15.dat"
$cfg.Range("B6").Value = "/home/ohel/Schreibtisch/FAST_sim/sim/NRELOffshrBsline5MW_InflowWind_12mps.dat"
$cfg.Range("B7").Value = "/home/ohel/Schreibtisch/FAST_sim/wind/TurbSim.inp"
$cfg.Range("B8").Value = "/home/ohel/Schreibtisch/FAST_sim/wind/IEC_template.IPT"
$cfg.Range("B10").Value = "/home/ohel/Schreibtisch/FAST_sim/sim"
$cfg.Range("B11").Value = "/home/ohel/Schreibtisch/FAST_sim/wind"

# Row 3 height change
$cfg.Rows.Item(3).RowHeight = 14.9

# Selection change
$cfg.Range("B29").Select()

# --- Sheet "DLC_List": update DLC values ---
$dlc = $wb.Worksheets.Item("DLC_List")

# Row 2: clear E2/F2 entirely (cell + formatting, not just the value)
$dlc.Range("E2").Clear()
$dlc.Range("F2").Clear()

# Row 4: E4 12 -> 8 ; F4 8 -> 3
$dlc.Range("E4").Value = "8"
$dlc.Range("F4").Value = "3"

# Row 5: clear E5/F5/H5 entirely (cell + formatting, not just the value)
$dlc.Range("E5").Clear()
$dlc.Range("F5").Clear()
$dlc.Range("H5").Clear()

# Column J width (target raw width 1.1; the engine snaps ColumnWidth to
# 1/7-character pixel increments, so 3/7 -> raw 8/7 = 1.142857... is the
# closest achievable approximation to 1.1)
$dlc.Columns.Item(10).ColumnWidth = 0.42857142857142855

# Selection change
$dlc.Range("F18").Select()
